$d = $word.ActiveDocument

# The opening "ID" marker paragraph currently reads
#   **ID__AFFARS_5309_topic_15__ID**<space>
# split across two runs (the second run holding just the trailing space).
# Replace it with the new marker text and drop the trailing space/run by
# including it in the search text but not in the replacement.
$d.Content.Find.Execute( `
    "**ID__AFFARS_5309_topic_15__ID** ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "**ID__AFFARS_5309_406_3__ID**", 2)

# Give that paragraph a light box border (5pt padding on every side) and
# bump its left indent from 6pt (120 twips) to 11.25pt (225 twips).
$p = $d.Paragraphs(1)
$pf = $p.Range.ParagraphFormat
$pf.LeftIndent = 11.25
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
